# Rebuild the "ESTADO DE CUENTA" worker/period table (rows 16-50).
# The table is regrouped by worker (C/D columns) instead of by period,
# several stale period rows are dropped, new more-recent periods are
# added, and GUILLERMO RODRIGUEZ PIÑERES's Salario Basico is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(16, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2401', 46400, 1160000),
    @(17, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2312', 46400, 1160000),
    @(18, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2311', 46400, 1160000),
    @(19, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2310', 46400, 1160000),
    @(20, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2309', 46400, 1160000),
    @(21, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2308', 46400, 1160000),
    @(22, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2307', 46400, 1160000),
    @(23, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2306', 46400, 1160000),
    @(24, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2305', 46400, 1160000),
    @(25, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2304', 46400, 1160000),
    @(26, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2303', 46400, 1160000),
    @(27, '9315662', 'ERICH RAFAEL HERRERA CABALLERO', '2302', 34027, 1160000),
    @(28, '1193522083', 'ESTIBENSON RAFAEL PALACIN VILLAREAL', '2211', 28000, 1000000),
    @(29, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2401', 46400, 1160000),
    @(30, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2312', 46400, 1160000),
    @(31, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2311', 46400, 1160000),
    @(32, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2310', 46400, 1160000),
    @(33, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2309', 46400, 1160000),
    @(34, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2308', 46400, 1160000),
    @(35, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2307', 46400, 1160000),
    @(36, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2306', 46400, 1160000),
    @(37, '1052219882', 'CARLOSMARIO SEPULVEDA PEDROZO', '2305', 38667, 1160000),
    @(38, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2401', 46400, 1200000),
    @(39, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2312', 46400, 1200000),
    @(40, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2311', 46400, 1200000),
    @(41, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2310', 46400, 1200000),
    @(42, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2309', 46400, 1200000),
    @(43, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2308', 46400, 1200000),
    @(44, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2307', 46400, 1200000),
    @(45, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2306', 46400, 1200000),
    @(46, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2305', 46400, 1200000),
    @(47, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2304', 46400, 1200000),
    @(48, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2303', 46400, 1200000),
    @(49, '73119562', 'GUILLERMO RODRIGUEZ PIÑERES', '2302', 15467, 1200000),
    @(50, '1143412227', 'LEIDY MARIAM BELLIDO ARROYO', '2303', 1547, 1160000),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value2  = $row[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value2  = $row[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value2  = $row[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value2  = $row[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value2  = $row[5]   # G - Salario Basico
}
